$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look like numbers,
# so they stay text (matching the original inlineStr string cells) instead of
# being auto-converted to numeric values by Excel.
$numericLookingCells = @("D4", "D5", "D6", "D8", "D9", "D11", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin data (prices, volumes, and the three reordered rows).
$ws.Range("D2").Value = '76.611.33'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '2.964.23'
$ws.Range("E3").Value = '  +3.25%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '198.89'
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").Value = '596.80'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.553'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").Value = '0.206'
$ws.Range("E9").Value = '  +7.76%  '
$ws.Range("D10").Value = '2.953.32'
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("D11").Value = '0.445'
$ws.Range("E11").Value = '  +11.55%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").Value = '4.92'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.496.06'
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = '28.60'
$ws.Range("E15").Value = '  +4.90%  '
$ws.Range("D16").Value = '76.464.64'
$ws.Range("E16").Value = '  +1.01%  '
$ws.Range("D17").Value = '0.0000192'
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").Value = '2.949.57'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '13.64'
$ws.Range("E19").Value = '  +8.93%  '
$ws.Range("D20").Value = '8.76'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("D21").Value = '380.08'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '2.30'
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").Value = '4.34'
$ws.Range("E23").Value = '  +4.99%  '
$ws.Range("D24").Value = '72.45'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.114.87'
$ws.Range("E25").Value = '  +2.36%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '4.30'
$ws.Range("E27").Value = '  +1.98%  '
$ws.Range("D28").Value = '9.76'
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").Value = '0.0000109'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '8.46'
$ws.Range("E31").Value = '  +9.19%  '
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").Value = '497.90'
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '20.34'
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '164.55'
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("D38").Value = '0.394'
$ws.Range("E38").Value = '  +14.77%  '
$ws.Range("D39").Value = '0.109'
$ws.Range("E39").Value = '  +20.03%  '
$ws.Range("D40").Value = '19.98'
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D43").Value = '180.53'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = '4.95'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '1.66'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").Value = '39.99'
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '0.595'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = '3.92'
$ws.Range("E49").Value = '  +4.14%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").Value = '0.668'
$ws.Range("E51").Value = '  +0.21%  '
